# The data table gained a new weekly record. Insert a new row at 513
# (this shifts the existing rows 513-602 down to 514-603, carrying their
# formatting/values with them) and populate the newly inserted row with
# the new observation's data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(513).Insert()

$ws.Range("A513").Value = 8
$ws.Range("B513").Value = "Terminal La Palmera de La Serena"
$ws.Range("C513").Value = "Coquimbo"
$ws.Range("D513").Value = 45180
$ws.Range("E513").Value = 4
$ws.Range("F513").Value = 100114013
$ws.Range("G513").Value = "Zanahoria"
$ws.Range("H513").Value = "Sin especificar"
$ws.Range("I513").Value = "Primera"
$ws.Range("J513").Value = 480
$ws.Range("K513").Value = 5500
$ws.Range("L513").Value = 6000
$ws.Range("M513").Value = 5750
$ws.Range("N513").Value = "$/saco 20 kilos"
$ws.Range("O513").Value = "Provincia del Elquí"
$ws.Range("P513").Value = 288
$ws.Range("Q513").Value = 20
$ws.Range("R513").Value = "Hortaliza"
